$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the previous data row (21) down to the new row (22)
# so the new row matches the date / number styles used throughout the table.
$ws.Range("A21:I21").Copy()
$ws.Range("A22:I22").PasteSpecial(-4122) # xlPasteFormats

# Append the new exercise record as row 22
$ws.Range("A22").Value = 44007
$ws.Range("B22").Value = 99
$ws.Range("C22").Value = 107
$ws.Range("D22").Value = 0.9
$ws.Range("E22").Value = "SAME"
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 82.4
$ws.Range("H22").Value = 28.9
$ws.Range("I22").Value = "OVERWEIGHT"

# Restore the active selection to C5, matching the saved workbook state
$ws.Range("C5").Select()
